# Swap the presentation's theme color scheme from the custom "Integral" /
# "Red Violet" palette over to the stock PowerPoint "Office Theme" / "Office"
# palette (the font scheme and format scheme are already identical between
# the two themes, so only the 12 theme colors need to change).

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$tcs = $theme.ThemeColorScheme

# Theme color scheme slot order (1-based):
#   1 dk1   2 lt1   3 dk2   4 lt2
#   5 accent1  6 accent2  7 accent3  8 accent4  9 accent5  10 accent6
#   11 hlink   12 folHlink
$tcs.Item(1).RGB  = 0x00 + (0x00 * 256) + (0x00 * 65536)   # dk1      000000
$tcs.Item(2).RGB  = 0xFF + (0xFF * 256) + (0xFF * 65536)   # lt1      FFFFFF
$tcs.Item(3).RGB  = 0x44 + (0x54 * 256) + (0x6A * 65536)   # dk2      44546A
$tcs.Item(4).RGB  = 0xE7 + (0xE6 * 256) + (0xE6 * 65536)   # lt2      E7E6E6
$tcs.Item(5).RGB  = 0x5B + (0x9B * 256) + (0xD5 * 65536)   # accent1  5B9BD5
$tcs.Item(6).RGB  = 0xED + (0x7D * 256) + (0x31 * 65536)   # accent2  ED7D31
$tcs.Item(7).RGB  = 0xA5 + (0xA5 * 256) + (0xA5 * 65536)   # accent3  A5A5A5
$tcs.Item(8).RGB  = 0xFF + (0xC0 * 256) + (0x00 * 65536)   # accent4  FFC000
$tcs.Item(9).RGB  = 0x44 + (0x72 * 256) + (0xC4 * 65536)   # accent5  4472C4
$tcs.Item(10).RGB = 0x70 + (0xAD * 256) + (0x47 * 65536)   # accent6  70AD47
$tcs.Item(11).RGB = 0x05 + (0x63 * 256) + (0xC1 * 65536)   # hlink    0563C1
$tcs.Item(12).RGB = 0x95 + (0x4F * 256) + (0x72 * 65536)   # folHlink 954F72
